$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two new "losers" columns (bold, like the rest of row 1)
$ws.Range("F1").Value = "2 losers"
$ws.Range("G1").Value = "3 losers"
$ws.Range("F1").Font.Bold = $true
$ws.Range("G1").Font.Bold = $true

# Old ELO values bumped from 100 to 1000 for players A-D
$ws.Range("B2").Value = 1000
$ws.Range("B3").Value = 1000
$ws.Range("B4").Value = 1000
$ws.Range("B5").Value = 1000

# Players E, F, G no longer carry an explicit old-ELO value
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()

# New ELO-delta formula: multiplier 8 -> 40, and (B-C) -> (C-B).
# Set one cell at a time (via R1C1) so each row keeps its own independent
# <f> instead of Excel collapsing the identical pattern into a shared formula.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 5).FormulaR1C1 = "=40*(RC[-1]-1/(1+10^((RC[-2]-RC[-3])/400)))"
}

# New columns: points lost when splitting the loss among 2 or 3 losers
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 6).FormulaR1C1 = "=RC[-1]/2"
    $ws.Cells.Item($r, 7).FormulaR1C1 = "=RC[-2]/3"
}

# Recalculate and refresh the selection / used range
$wb.Application.Calculate()
$ws.Range("F10").Select()
